# Rebuild the deck: a title (question) slide followed by six
# "Title and Content" slides carrying the strategy-report outline.

$p = $ppt.ActivePresentation

# --- Slide 1: replace the old "Title and Content" slide with a
#     "Title Slide" (ctrTitle + subTitle) layout carrying the new prompt.
$p.Slides.Item(1).Delete()
$s1 = $p.Slides.Add(1, 1)
$s1.Shapes.Item(1).TextFrame.TextRange.Text = "在筆電生產效率遇到瓶頸且人工組裝速度較慢的情境下，如何通過導入AI技術來提升生產效率10%？"
# Subtitle stays empty, matching the target deck (leave the placeholder's
# default empty paragraph untouched).

# --- Slides 2-7: "Title and Content" slides.
$defs = @(
    @{
        Title = "目標"
        Bullets = @(
            "主動識別並解決人工組裝過程中的不效率問題。",
            "利用AI技術優化生產流程，並實現可持續的效率提升。",
            "促進跨部門合作，確保技術導入的成功。"
        )
    },
    @{
        Title = "梯形分析：尋找問題甜蜜點"
        Bullets = @(
            "我們如何提升筆電的整體生產效率？",
            "我們如何提升生產流程中每個組裝工序的效率？",
            "我們如何通過自動化與AI技術來替代或加速人工組裝？"
        )
    },
    @{
        Title = "痛點"
        Bullets = @(
            "人工組裝速度慢，易造成生產延遲。",
            "現有生產流程缺乏數據支持，導致難以進行效率分析與優化。",
            "對於AI技術的認知不足，許多員工對於新技術的抗拒態度。"
        )
    },
    @{
        Title = "跨部門視角的整合分析"
        Bullets = @(
            "總經理：導入AI技術將提升整體生產力，替公司創造競爭優勢，並提高市場份額。",
            "生產部門：需針對現有流程進行數據收集，分析瓶頸之處，提供合適的數據支持。",
            "IT部門：需掌握開發機械手臂的技術，並聯繫外部資源以獲取AI系統的開發支持。",
            "人力資源部門：需針對員工進行必要的AI技術和自動化操作培訓，以降低對新技術的抵抗。"
        )
    },
    @{
        Title = "實作步驟"
        Bullets = @(
            "需求分析：各部門會議，確定AI導入的具體需求與目標。",
            "數據收集：生產部門提供現有流程的數據，IT部門分析數據以找出關鍵問題。",
            "技術開發：IT部門與AI公司合作，開發機械手臂及AI系統，針對特定組裝工序進行優化。",
            "員工培訓：人力資源部門組織培訓，引導員工熟悉新技術與工具的使用。",
            "階段性評估：在導入過程中持續監測進展，根據數據分析及時調整策略以達成預期的效率提升。"
        )
    },
    @{
        Title = "結論"
        Bullets = @(
            "通過以上策略和實作步驟，筆電的生產效率有望在導入AI技術後提升10%。",
            "各部門的協同合作是成功的關鍵，特別是生產部門的數據支持和IT部門的技術開發。",
            "確保員工的適應性與接受度將直接影響AI實施的成效，因此必須重視培訓與支持。"
        )
    }
)

$idx = 2
foreach ($d in $defs) {
    $slide = $p.Slides.Add($idx, 2)
    $slide.Shapes.Item(1).TextFrame.TextRange.Text = $d.Title
    $body = $slide.Shapes.Item(2).TextFrame.TextRange
    $body.Text = "`r" + ($d.Bullets -join "`r")
    $idx = $idx + 1
}

Write-Host "Slides:" $p.Slides.Count
